$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bhaskar Lalwani
$ws.Range("C2").Value = 82.2
$ws.Range("D2").Value = 46
$ws.Range("E2").Value = 37

# Row 3 - Mayur Gogoi
$ws.Range("C3").Value = 74.09999999999999
$ws.Range("D3").Value = 47
$ws.Range("E3").Value = 34

# Row 4 - Aniruddha Mukherjee
$ws.Range("C4").Value = 88.09999999999999
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 44

# Row 5 - Amandeep Chourasia
$ws.Range("C5").Value = 69.90000000000001
$ws.Range("D5").Value = 48
$ws.Range("E5").Value = 33

# Row 6 - Ishaan Mukherjee
$ws.Range("C6").Value = 84.40000000000001
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 33
